# Commit: "User can follow sites, Ping : progress"
#
# The "sites" and "openings" sheets are restructured from small lookup
# tables (name / short-code / numeric-code / extra columns) into plain
# single-column lists of values (sites gains a couple of new entries:
# 046P / Pääkeittiö; openings gains 046G-S, and kontti becomes Kontti).
# "users" / "groups" keep their data untouched - they only shift which
# shared-string index they point at, which happens automatically once
# the now-unused strings (pääkeittiö, kontti, Kansi 3, Kansi 4, SITE,
# OPENING, 057G) are dropped from the shared-string table as a side
# effect of rewriting the other two sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# sites
# ---------------------------------------------------------------
$sites = $wb.Worksheets.Item("sites")

# Drop the old 8 rows (A:D) entirely - Rows(...).Delete() shifts the
# remaining rows up without disturbing the sheet's <cols> layout.
$sites.Rows("1:8").Delete() | Out-Null

$sites.Range("A1").Value = "K3"
$sites.Range("A2").Value = "036G"
$sites.Range("A3").Value = "36G"
$sites.Range("A4").Value = 36
$sites.Range("A5").Value = "K4"
$sites.Range("A6").Value = "046P"
$sites.Range("A7").Value = "046G"
$sites.Range("A8").Value = "46G"
$sites.Range("A9").Value = 46
$sites.Range("A10").Value = "047G"
$sites.Range("A11").Value = "47G"
$sites.Range("A12").Value = 47
$sites.Range("A13").Value = "Pääkeittiö"
$sites.Range("A14").Value = "S"

# ---------------------------------------------------------------
# openings
# ---------------------------------------------------------------
$openings = $wb.Worksheets.Item("openings")

$openings.Rows("1:5").Delete() | Out-Null

$openings.Range("A1").Value = "035N-S"
$openings.Range("A2").Value = 35
$openings.Range("A3").Value = "036N-S"
$openings.Range("A4").Value = 36
$openings.Range("A6").Value = 46
$openings.Range("A7").Value = "047G-P"
$openings.Range("A8").Value = 47
$openings.Range("A9").Value = "Kontti"
$openings.Range("A5").Value = "046G-S"
$openings.Range("A10").Value = "O"

# ---------------------------------------------------------------
# Selection / view state (mirrors the new activeCell/sqref in the diff)
# ---------------------------------------------------------------
$sites.Range("D10").Select() | Out-Null
$openings.Range("D12").Select() | Out-Null

# Restore the workbook's active tab back to "users" (xr:revisionPtr /
# workbookView activeTab="2" in the diff keeps pointing at "users" -
# Range.Select() above activates whichever sheet it runs on, so re-pick
# the users sheet last to leave it as the active / tabSelected one).
$users = $wb.Worksheets.Item("users")
$users.Activate() | Out-Null
$users.Range("A2").Select() | Out-Null
